# edit.ps1 -- apply the CV update described by the commit:
#   1. "React.js and AWS" -> "React and AWS" (bullet under the Lab714 entry)
#   2. hyperlink display text "linkedin.com/in/aadarsha2002/" -> "aadarsha2002.github.io"
#      (keep the run's existing formatting: color 0563C1 + single underline)
#
# (A third part of the upstream commit bundles a Word "Insert Add-in" task-pane
# reference (word/webextensions/...). That is a raw OOXML package part that is
# not reachable from the Word object model / COM automation surface -- there is
# no Application/Document method that mints new parts like that (TaskPanes and
# CustomXMLParts are read-only collections here, same as in real Word COM) --
# so it is intentionally left out of this script.)

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Developing software using React.js and AWS ..." -> "...React and AWS ..."
# ---------------------------------------------------------------------------
$oldSkills = "Developing software using React.js and AWS that works with proprietary IoT devices to extract, organize, and analyze data."
$newSkills = "Developing software using React and AWS that works with proprietary IoT devices to extract, organize, and analyze data."

foreach ($story in $d.StoryRanges) {
    $probe = $story.Duplicate
    $found = $probe.Find.Execute($oldSkills, $true, $false, $false, $false, $false, $true, 1, $false, $newSkills, 2)
    if ($found) {
        Write-Output ("Replaced React.js wording in story type " + $story.StoryType)
    }
}

# ---------------------------------------------------------------------------
# 2) hyperlink display text in the contact-info header:
#    "linkedin.com/in/aadarsha2002/" -> "aadarsha2002.github.io"
#    Done via explicit Range bounds (rather than letting Find itself perform
#    the replace) so we fully control the resulting run; the original run's
#    look (blue 0563C1, single underline) is re-applied explicitly afterward
#    to match what the source document already showed for that hyperlink.
# ---------------------------------------------------------------------------
$oldUrl = "linkedin.com/in/aadarsha2002/"
$newUrl = "aadarsha2002.github.io"
$linkColor = 12673797   # RGB(5,99,193) == hex 0563C1, stored BGR for WdColor
$linkUnderline = 1      # wdUnderlineSingle

foreach ($story in $d.StoryRanges) {
    $storyText = $story.Text
    $pos = $storyText.IndexOf($oldUrl)
    if ($pos -ge 0) {
        $startChar = $story.Start + $pos
        $endChar = $startChar + $oldUrl.Length

        $target = $story.Duplicate
        $target.Start = $startChar
        $target.End = $endChar

        $target.Text = $newUrl
        $target.Font.Color = $linkColor
        $target.Font.Underline = $linkUnderline

        Write-Output ("Replaced hyperlink text in story type " + $story.StoryType + ": " + $target.Text)
    }
}
